$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-5 (Neutrophils / Resolving-Mac combinations no longer present)
$ws.Range("A3:T5").Delete()

# Update the recomputed TPM-based values in row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.195179
$ws.Range("N2").Value = 0.585537
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.198928974127
$ws.Range("R2").Value = 1.790360767143
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
